# "Added Test Data Model" - add an "Expected Result" column (C) to the
# HomePage test-data table, mirroring the existing Selenium/Cypress values
# already present in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Expected Result"
$ws.Range("C2").Value = "Selenium"
$ws.Range("C3").Value = "Cypress"

# Match the look of the existing table: copy the header/body formatting
# from column B (bold header with fill/border, bordered body cells) onto
# the new column C.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null

$ws.Range("B2:B3").Copy() | Out-Null
$ws.Range("C2:C3").PasteSpecial(-4122) | Out-Null

# Size the new column to fit its contents (header "Expected Result" -> 14).
$ws.Columns.Item(3).ColumnWidth = 13.1666666666667

$ws.Range("I9").Select() | Out-Null
